$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.338.19"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.711.23"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5290"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06679"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2660"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.81"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07701"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.504"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.946.41"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.710.59"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5853"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "0.0₅8212"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.97"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "27.361.07"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "223.01"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.48"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.010"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.45"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.689"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1207"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05338"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.464"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.430"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.636"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.873"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9521"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.395"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("D39").Value = "1.148.61"
$ws.Range("E39").Value = "  +8.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01635"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.788"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8396"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "1.853.59"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4563"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.109"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05215"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.44%  "
